$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.683.93'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '3.414.10'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '569.82'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").Value = '157.38'
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.412.68'
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("E11").Value = '  +4.46%  '
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '3.999.65'
$ws.Range("E13").Value = '  +2.28%  '
$ws.Range("E14").Value = '  -2.99%  '
$ws.Range("E15").Value = '  +9.09%  '
$ws.Range("D16").Value = '27.26'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '63.684.91'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").Value = '3.392.92'
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").Value = '6.27'
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").Value = '378.20'
$ws.Range("E21").Value = '  -1.19%  '
$ws.Range("E22").Value = '  -3.34%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '71.92'
$ws.Range("E24").Value = '  +2.95%  '
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("E26").Value = '  +28.81%  '
$ws.Range("D27").Value = '9.41'
$ws.Range("E27").Value = '  +4.91%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = '6.06'
$ws.Range("E30").Value = '  +8.55%  '
$ws.Range("D31").Value = '1.37'
$ws.Range("E31").Value = '  +5.47%  '
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '6.35'
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '6.79'
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").Value = '158.74'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("D39").Value = '2.978.36'
$ws.Range("E39").Value = '  +6.85%  '
$ws.Range("E40").Value = '  +3.14%  '
$ws.Range("D41").Value = '27.00'
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("D44").Value = '41.95'
$ws.Range("E44").Value = '  +3.92%  '
$ws.Range("D45").Value = '0.761'
$ws.Range("E45").Value = '  +2.75%  '
$ws.Range("D46").Value = '4.32'
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("E47").Value = '  +5.84%  '
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  +3.66%  '
$ws.Range("E49").Value = '  +23.85%  '
$ws.Range("D50").Value = '295.07'
$ws.Range("E50").Value = '  +3.03%  '
$ws.Range("D51").Value = '6.34'
$ws.Range("E51").Value = '  +0.99%  '
